$wb = $excel.ActiveWorkbook

$oldGuid = "4d9f9193-0a05-4f8d-8d4a-4af7baa7b82b"
$newGuid = "91171a96-5f7d-4641-8576-44b0058bbe25"

$oldHoDate = "2016-08-29 02:57:06"
$newHoDate = "2016-08-29 02:57:33"

$oldZhHash = "845ac07b330cdf7371dd087b0eb335c410cc64d3"
$newZhHash = "9800bf1d99ddd4eb7800e02a0ef556eabf3e0ef6"

$oldZhDate = "2016-08-29 02:56:57"
$newZhDate = "2016-08-29 02:57:28"

$oldDeHash = "845ac07b330cdf7371dd087b0eb335c410cc64d3"
$newDeHash = "9800bf1d99ddd4eb7800e02a0ef556eabf3e0ef6"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1906dd84a97b5ec0f7db97ab02bb69836d17f9b3/e2e/$oldGuid.md"

# ---------- Overview sheet ----------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = $newHoDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------- zh-cn sheet ----------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$ws2.Range("H2").Value = $newZhDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md") | Out-Null

# ---------- de-de sheet ----------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$ws3.Range("H2").Value = $newHoDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md") | Out-Null
